$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 1559
$ws.Range("F2").Value = 37364696
$ws.Range("G2").Value = 95915927

$ws.Range("E3").Value = 1599
$ws.Range("F3").Value = 38755392
$ws.Range("G3").Value = 47251260

$ws.Range("E4").Value = 1532
$ws.Range("F4").Value = 38637944
$ws.Range("G4").Value = 33316157

$ws.Range("E5").Value = 1578
$ws.Range("F5").Value = 39457520
$ws.Range("G5").Value = 33780466

$ws.Range("E6").Value = 1626
$ws.Range("F6").Value = 39015272
$ws.Range("G6").Value = 24163686

$ws.Range("D9").Select()
